$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update unit prices (column G) ---
$ws.Range("G2").Value = 0.033
$ws.Range("G3").Value = 0.03
$ws.Range("G4").Value = 0.31
$ws.Range("G5").Value = 0.269
$ws.Range("G6").Value = 0.103
$ws.Range("G7").Value = 0.062
$ws.Range("G8").Value = 0.135
$ws.Range("G9").Value = 0.009
$ws.Range("G10").Value = 0.012

# --- Remove the existing hyperlink on I4 (it will be re-added on I5) ---
$ws.Range("I4").Hyperlinks.Delete()

# --- Populate the new datasheet / supplier links in column I ---
# Order matters: it controls how new shared strings are appended.
$ws.Range("I4").Value = "https://au.element14.com/panasonic/eeefc1c331p/cap-330-f-16v-radial-smd/dp/9694420?st=330%20micro%20farad%20capacitor"
$ws.Range("I3").Value = "https://au.element14.com/multicomp/mc0805b331k500ct/cap-330pf-50v-10-x7r-0805/dp/1759212"
$ws.Range("I2").Value = "https://au.element14.com/avx/08055c104kat2a/cap-0-1-f-50v-10-x7r-0805-reel/dp/2280990"

$ws.Range("I5").Value = "https://au.element14.com/tdk/c3216x7r1h105k160ab/cap-1-f-50v-10-x7r-1206/dp/1907358"
$ws.Hyperlinks.Add($ws.Range("I5"), "https://au.mouser.com/ProductDetail/Nichicon/UWT1E331MNL1GS?qs=sGAEpiMZZMsh%252B1woXyUXjxib3yTgZFe2bgHsEDyVhdM%3D")
$ws.Range("I5").Style = "Hyperlink"

$ws.Range("I6").Value = "https://au.element14.com/vishay/crcw080510k0fkea/res-10k-1-0-125w-0805-thick-film/dp/1469856"
$ws.Range("I7").Value = "https://au.element14.com/yageo/rc0805fr-073k3l/res-3k3-1-0-125w-0805-thick-film/dp/9237682"
$ws.Range("I8").Value = "https://au.element14.com/vishay/crcw080539k0fkeahp/res-39k-1-0-33w-0805-thick-film/dp/1738980"
$ws.Range("I9").Value = "https://au.element14.com/multicomp/mc01w080551k/res-1k-5-0-1w-0805-thick-film/dp/9333711"
$ws.Range("I10").Value = "https://au.element14.com/multicomp/mcmr08x3300ftl/res-330r-1-0-125w-0805-ceramic/dp/2073741"

# --- Update the selected range shown when the workbook is opened ---
$ws.Range("H2:H13").Select()
